$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price (D) and volume (E) columns are treated as plain text so that
# values like "130.60" or "2.903.04" are preserved exactly, without Excel
# re-interpreting them as numbers (which would drop trailing zeros / add float noise).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.837.94'
$ws.Range("E2").Value = '  -3.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.907.09'
$ws.Range("E3").Value = '  -3.90%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.07'
$ws.Range("E5").Value = '  -0.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.80'
$ws.Range("E6").Value = '  -5.33%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.505'
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.906.55'
$ws.Range("E9").Value = '  -3.81%  '
$ws.Range("E10").Value = '  -4.05%  '
$ws.Range("E11").Value = '  -4.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.444'
$ws.Range("E12").Value = '  -4.03%  '
$ws.Range("E13").Value = '  -2.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.52'
$ws.Range("E14").Value = '  -5.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.127'
$ws.Range("E15").Value = '  +1.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.383.42'
$ws.Range("E16").Value = '  -4.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.763.95'
$ws.Range("E17").Value = '  -3.45%  '
$ws.Range("E18").Value = '  -4.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.901.25'
$ws.Range("E19").Value = '  -4.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '430.48'
$ws.Range("E20").Value = '  -4.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.53'
$ws.Range("E21").Value = '  -4.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.683'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.10'
$ws.Range("E23").Value = '  -5.71%  '
$ws.Range("E24").Value = '  -1.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.87'
$ws.Range("E25").Value = '  -4.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.22'
$ws.Range("E26").Value = '  -3.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.03'
$ws.Range("E27").Value = '  -2.56%  '
$ws.Range("E29").Value = '  +0.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("E31").Value = '  -2.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.07'
$ws.Range("E32").Value = '  -5.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.57'
$ws.Range("E33").Value = '  -3.97%  '
$ws.Range("E34").Value = '  -2.90%  '
$ws.Range("E35").Value = '  -2.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("E36").Value = '  -3.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.61'
$ws.Range("E37").Value = '  -4.61%  '
$ws.Range("E38").Value = '  -4.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.55'
$ws.Range("E39").Value = '  -1.89%  '
$ws.Range("E40").Value = '  -4.31%  '
$ws.Range("E41").Value = '  -4.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.63'
$ws.Range("E42").Value = '  -4.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.292'
$ws.Range("E43").Value = '  -4.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.95'
$ws.Range("E44").Value = '  -10.39%  '
$ws.Range("E45").Value = '  -3.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '373.41'
$ws.Range("E46").Value = '  -4.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.703.52'
$ws.Range("E47").Value = '  -0.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.60'
$ws.Range("E48").Value = '  -2.44%  '
$ws.Range("E50").Value = '  -10.22%  '
$ws.Range("E51").Value = '  -2.16%  '
